$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting rows 116:172 down to 117:173
$ws.Rows.Item(116).Insert()

# Fill the new row 116 with data
$ws.Cells.Item(116, 1).Value = 7
$ws.Cells.Item(116, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(116, 3).Value = "Ñuble"
$ws.Cells.Item(116, 4).Value = 44460
$ws.Cells.Item(116, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat
$ws.Cells.Item(116, 5).Value = 16
$ws.Cells.Item(116, 6).Value = 100114013
$ws.Cells.Item(116, 7).Value = "Zanahoria"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 160
$ws.Cells.Item(116, 11).Value = 8500
$ws.Cells.Item(116, 12).Value = 9000
$ws.Cells.Item(116, 13).Value = 8750
$ws.Cells.Item(116, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(116, 16).Value = 438
$ws.Cells.Item(116, 17).Value = 20
$ws.Cells.Item(116, 18).Value = "Hortaliza"
